$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 150
$ws.Range("I9").Value = 150
$ws.Range("K9").Value = 150
$ws.Range("M9").Value = 19
$ws.Range("H98").Value = 1290.7407
$ws.Range("I98").Value = 1314.04
$ws.Range("K98").Value = 1314.04
$ws.Range("M98").Value = 183.96
$ws.Range("H116").Value = 4942.25
$ws.Range("I116").Value = 4941
$ws.Range("K116").Value = 4941
$ws.Range("M116").Value = -1499
$ws.Range("H122").Value = 1290.7407
$ws.Range("I122").Value = 1314.04
$ws.Range("K122").Value = 3942.12
$ws.Range("M122").Value = -1492.12
$ws.Range("H132").Value = 1542.125
$ws.Range("I132").Value = 1676.8572
$ws.Range("K132").Value = 5030.571599999999
$ws.Range("M132").Value = -2500.571599999999
$ws.Range("H138").Value = 2506.611
$ws.Range("J138").Value = 3880
$ws.Range("L138").Value = 11640
$ws.Range("N138").Value = -21920

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3042.5715
$ws.Range("I45").Value = 2474.5
$ws.Range("K45").Value = 2474.5
$ws.Range("M45").Value = -2097.5
$ws.Range("H74").Value = 528.03705
$ws.Range("I74").Value = 528.03705
$ws.Range("K74").Value = 528.03705
$ws.Range("M74").Value = 345.96295
$ws.Range("H77").Value = 528.03705
$ws.Range("I77").Value = 528.03705
$ws.Range("K77").Value = 2640.18525
$ws.Range("M77").Value = 1727.81475
$ws.Range("H97").Value = 517.7826
$ws.Range("I97").Value = 492.2381
$ws.Range("J97").Value = 786
$ws.Range("K97").Value = 492.2381
$ws.Range("L97").Value = 786
$ws.Range("M97").Value = 3.761900000000026
$ws.Range("N97").Value = -1778
$ws.Range("H101").Value = 25734.666
$ws.Range("J101").Value = 25734.666
$ws.Range("L101").Value = 25734.666
$ws.Range("N101").Value = -32224.666
$ws.Range("H122").Value = 1140373.4
$ws.Range("I122").Value = 1457980
$ws.Range("J122").Value = 28750
$ws.Range("K122").Value = 4373940
$ws.Range("L122").Value = 86250
$ws.Range("M122").Value = -4371490
$ws.Range("N122").Value = -91150

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 50000
$ws.Range("I60").Value = 30000
$ws.Range("J60").Value = 60000
$ws.Range("K60").Value = 30000
$ws.Range("L60").Value = 60000
$ws.Range("M60").Value = -29401
$ws.Range("N60").Value = -61198
$ws.Range("H86").Value = 3068.3333
$ws.Range("I86").Value = 3082
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 3082
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -1959
$ws.Range("N86").Value = -5246
$ws.Range("H89").Value = 3068.3333
$ws.Range("I89").Value = 3082
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 15410
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -9794
$ws.Range("N89").Value = -26232
$ws.Range("H100").Value = 19910.5
$ws.Range("J100").Value = 19910.5
$ws.Range("L100").Value = 19910.5
$ws.Range("N100").Value = -22074.5
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3784.261
$ws.Range("I31").Value = 2532.8235
$ws.Range("K31").Value = 2532.8235
$ws.Range("M31").Value = -2237.8235
$ws.Range("H34").Value = 3784.261
$ws.Range("I34").Value = 2532.8235
$ws.Range("K34").Value = 2532.8235
$ws.Range("M34").Value = -2330.8235
$ws.Range("H58").Value = 2607.2173
$ws.Range("I58").Value = 1374.875
$ws.Range("K58").Value = 1374.875
$ws.Range("M58").Value = -1171.875
$ws.Range("H74").Value = 46712.285
$ws.Range("J74").Value = 46712.285
$ws.Range("L74").Value = 46712.285
$ws.Range("N74").Value = -48460.285
$ws.Range("H77").Value = 46712.285
$ws.Range("J77").Value = 46712.285
$ws.Range("L77").Value = 140136.855
$ws.Range("N77").Value = -148872.855
$ws.Range("H134").Value = 2728.2
$ws.Range("I134").Value = 2402.3635
$ws.Range("J134").Value = 3624.25
$ws.Range("K134").Value = 7207.0905
$ws.Range("L134").Value = 10872.75
$ws.Range("M134").Value = -4672.0905
$ws.Range("N134").Value = -15942.75
$ws.Range("H136").Value = 2607.2173
$ws.Range("I136").Value = 1374.875
$ws.Range("K136").Value = 4124.625
$ws.Range("M136").Value = -1574.625

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 749
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 749
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 2247
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -2823
$ws.Range("H121").Value = 669.2857
$ws.Range("I121").Value = 697.6
$ws.Range("K121").Value = 2092.8
$ws.Range("M121").Value = -782.8000000000002

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3606.6924
$ws.Range("I80").Value = 2626.4285
$ws.Range("J80").Value = 4750.3335
$ws.Range("K80").Value = 2626.4285
$ws.Range("L80").Value = 4750.3335
$ws.Range("M80").Value = -1628.4285
$ws.Range("N80").Value = -6746.3335
$ws.Range("H83").Value = 3606.6924
$ws.Range("I83").Value = 2626.4285
$ws.Range("J83").Value = 4750.3335
$ws.Range("K83").Value = 13132.1425
$ws.Range("L83").Value = 23751.6675
$ws.Range("M83").Value = -8140.1425
$ws.Range("N83").Value = -33735.6675
$ws.Range("H98").Value = 9248.25
$ws.Range("J98").Value = 9248.25
$ws.Range("L98").Value = 9248.25
$ws.Range("N98").Value = -15238.25
$ws.Range("H113").Value = 67497.336
$ws.Range("I113").Value = 51247
$ws.Range("K113").Value = 51247
$ws.Range("M113").Value = -49077
$ws.Range("H122").Value = 251770.5
$ws.Range("I122").Value = 2361
$ws.Range("J122").Value = 999999
$ws.Range("K122").Value = 7083
$ws.Range("L122").Value = 2999997
$ws.Range("M122").Value = -4633
$ws.Range("N122").Value = -3004897
$ws.Range("H126").Value = 4999
$ws.Range("I126").Value = 4999
$ws.Range("K126").Value = 14997
$ws.Range("M126").Value = -12527
$ws.Range("H132").Value = 3333.3333
$ws.Range("I132").Value = 3000
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -17060

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 23199.8
$ws.Range("I47").Value = 22000
$ws.Range("J47").Value = 23499.75
$ws.Range("K47").Value = 22000
$ws.Range("L47").Value = 23499.75
$ws.Range("M47").Value = -21510
$ws.Range("N47").Value = -24479.75
$ws.Range("H52").Value = 23199.8
$ws.Range("I52").Value = 22000
$ws.Range("J52").Value = 23499.75
$ws.Range("K52").Value = 22000
$ws.Range("L52").Value = 23499.75
$ws.Range("M52").Value = -21767
$ws.Range("N52").Value = -23965.75
$ws.Range("H68").Value = 3100.25
$ws.Range("J68").Value = 4001.5
$ws.Range("L68").Value = 4001.5
$ws.Range("N68").Value = -5499.5
$ws.Range("H71").Value = 3100.25
$ws.Range("J71").Value = 4001.5
$ws.Range("L71").Value = 20007.5
$ws.Range("N71").Value = -27495.5
$ws.Range("H82").Value = 914.13336
$ws.Range("I82").Value = 776
$ws.Range("J82").Value = 1294
$ws.Range("K82").Value = 776
$ws.Range("L82").Value = 1294
$ws.Range("M82").Value = -415
$ws.Range("N82").Value = -2016
$ws.Range("H85").Value = 914.13336
$ws.Range("I85").Value = 776
$ws.Range("J85").Value = 1294
$ws.Range("K85").Value = 776
$ws.Range("L85").Value = 1294
$ws.Range("M85").Value = 472
$ws.Range("N85").Value = -3790
$ws.Range("H93").Value = 1443
$ws.Range("J93").Value = 1871.5
$ws.Range("L93").Value = 1871.5
$ws.Range("N93").Value = -4367.5
$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H48").Value = 44410
$ws.Range("J48").Value = 43012.5
$ws.Range("L48").Value = 43012.5
$ws.Range("N48").Value = -44150.5
$ws.Range("H68").Value = 59998
$ws.Range("J68").Value = 59998
$ws.Range("L68").Value = 59998
$ws.Range("N68").Value = -61620
$ws.Range("H71").Value = 59998
$ws.Range("J71").Value = 59998
$ws.Range("L71").Value = 179994
$ws.Range("N71").Value = -188106
$ws.Range("H113").Value = 711.8570999999999
$ws.Range("I113").Value = 598.5
$ws.Range("J113").Value = 863
$ws.Range("K113").Value = 1795.5
$ws.Range("L113").Value = 2589
$ws.Range("M113").Value = 374.5
$ws.Range("N113").Value = -6929
$ws.Range("H122").Value = 2674.75
$ws.Range("I122").Value = 2674.75
$ws.Range("K122").Value = 8024.25
$ws.Range("M122").Value = -5574.25
